$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Scenarios")

# Clear the "Number of Test Cases" values in column E for the data rows
# (E11:E41) — this also removes the stray "09" text value that lived in
# E18, which is why the shared string table shrinks by one entry.
$ws.Range("E11:E41").ClearContents()

# Update the sheet view: scroll back to the top, change zoom, and move
# the active selection.
$ws.Activate() | Out-Null
$window = $excel.ActiveWindow
$window.ScrollRow = 1
$window.ScrollColumn = 1
$window.Zoom = 84
$ws.Range("I19").Select() | Out-Null
